# Auto update Excel log
# Appends newly-logged sensor readings (2026-01-28, ~11:58-12:03) to the
# PIR, Humidity and Temperature sheets, extending each sheet's used range.
#
# Each data row is @(RowNumber, Date, Timestamp, Hour, Location, Value, Status)

$wb = $excel.ActiveWorkbook

$pirRows = @(
    @(146, '2026-01-28','11:58:34','11:00','Bathroom','No Motion','Inactive'),
    @(147, '2026-01-28','11:58:39','11:00','Bathroom','No Motion','Inactive'),
    @(148, '2026-01-28','11:58:44','11:00','Bathroom','No Motion','Inactive'),
    @(149, '2026-01-28','11:58:49','11:00','Bathroom','No Motion','Inactive'),
    @(150, '2026-01-28','11:58:54','11:00','Bathroom','No Motion','Inactive'),
    @(151, '2026-01-28','11:58:58','11:00','Bathroom','Motion Detected','Active'),
    @(152, '2026-01-28','12:02:50','12:00','Bathroom','No Motion','Inactive'),
    @(153, '2026-01-28','12:02:55','12:00','Bathroom','No Motion','Inactive'),
    @(154, '2026-01-28','12:03:00','12:00','Bathroom','No Motion','Inactive'),
    @(155, '2026-01-28','12:03:05','12:00','Bathroom','No Motion','Inactive'),
    @(156, '2026-01-28','12:03:10','12:00','Bathroom','No Motion','Inactive'),
    @(157, '2026-01-28','12:03:15','12:00','Bathroom','No Motion','Inactive'),
    @(158, '2026-01-28','12:03:20','12:00','Bathroom','No Motion','Inactive'),
    @(159, '2026-01-28','12:03:25','12:00','Bathroom','No Motion','Inactive'),
    @(160, '2026-01-28','12:03:30','12:00','Bathroom','No Motion','Inactive'),
    @(161, '2026-01-28','12:03:35','12:00','Bathroom','No Motion','Inactive'),
    @(162, '2026-01-28','12:03:40','12:00','Bathroom','No Motion','Inactive'),
    @(163, '2026-01-28','12:03:45','12:00','Bathroom','No Motion','Inactive')
)

$humidityRows = @(
    @(134, '2026-01-28','11:58:33','11:00','Bathroom','88.3%','Active'),
    @(135, '2026-01-28','11:58:37','11:00','Bathroom','88.3%','Active'),
    @(136, '2026-01-28','11:58:41','11:00','Bathroom','87.4%','Active'),
    @(137, '2026-01-28','11:58:45','11:00','Bathroom','88.3%','Active'),
    @(138, '2026-01-28','11:58:53','11:00','Bathroom','87.4%','Active'),
    @(139, '2026-01-28','11:58:57','11:00','Bathroom','88.4%','Active'),
    @(140, '2026-01-28','11:59:01','11:00','Bathroom','87.4%','Active'),
    @(141, '2026-01-28','12:02:50','12:00','Bathroom','88.0%','Active'),
    @(142, '2026-01-28','12:02:54','12:00','Bathroom','87.9%','Active'),
    @(143, '2026-01-28','12:03:02','12:00','Bathroom','87.9%','Active'),
    @(144, '2026-01-28','12:03:06','12:00','Bathroom','87.0%','Active'),
    @(145, '2026-01-28','12:03:10','12:00','Bathroom','87.9%','Active'),
    @(146, '2026-01-28','12:03:14','12:00','Bathroom','87.9%','Active'),
    @(147, '2026-01-28','12:03:18','12:00','Bathroom','87.9%','Active'),
    @(148, '2026-01-28','12:03:26','12:00','Bathroom','87.0%','Active'),
    @(149, '2026-01-28','12:03:30','12:00','Bathroom','88.0%','Active'),
    @(150, '2026-01-28','12:03:34','12:00','Bathroom','88.0%','Active'),
    @(151, '2026-01-28','12:03:38','12:00','Bathroom','87.0%','Active'),
    @(152, '2026-01-28','12:03:42','12:00','Bathroom','88.0%','Active'),
    @(153, '2026-01-28','12:03:47','12:00','Bathroom','87.1%','Active')
)

$temperatureRows = @(
    @(134, '2026-01-28','11:58:33','11:00','Bathroom','22.9C','Active'),
    @(135, '2026-01-28','11:58:37','11:00','Bathroom','22.9C','Active'),
    @(136, '2026-01-28','11:58:41','11:00','Bathroom','22.9C','Active'),
    @(137, '2026-01-28','11:58:45','11:00','Bathroom','22.9C','Active'),
    @(138, '2026-01-28','11:58:53','11:00','Bathroom','22.9C','Active'),
    @(139, '2026-01-28','11:58:57','11:00','Bathroom','22.9C','Active'),
    @(140, '2026-01-28','11:59:01','11:00','Bathroom','22.9C','Active'),
    @(141, '2026-01-28','12:02:51','12:00','Bathroom','23.0C','Active'),
    @(142, '2026-01-28','12:02:55','12:00','Bathroom','23.0C','Active'),
    @(143, '2026-01-28','12:03:03','12:00','Bathroom','23.0C','Active'),
    @(144, '2026-01-28','12:03:07','12:00','Bathroom','23.0C','Active'),
    @(145, '2026-01-28','12:03:11','12:00','Bathroom','23.0C','Active'),
    @(146, '2026-01-28','12:03:15','12:00','Bathroom','23.0C','Active'),
    @(147, '2026-01-28','12:03:19','12:00','Bathroom','23.0C','Active'),
    @(148, '2026-01-28','12:03:27','12:00','Bathroom','23.0C','Active'),
    @(149, '2026-01-28','12:03:31','12:00','Bathroom','23.0C','Active'),
    @(150, '2026-01-28','12:03:35','12:00','Bathroom','23.0C','Active'),
    @(151, '2026-01-28','12:03:39','12:00','Bathroom','23.0C','Active'),
    @(152, '2026-01-28','12:03:43','12:00','Bathroom','23.0C','Active'),
    @(153, '2026-01-28','12:03:47','12:00','Bathroom','23.0C','Active')
)

# Writes one log row to $ws at $entry[0]. $percentCol5 = $true means column E
# (Value) also needs a quote-prefix to stay text (e.g. Humidity's "88.3%"
# would otherwise be auto-converted to the number 0.883).
function Write-LogRow($ws, $entry, $percentCol5) {
    $r = $entry[0]

    # Column A (Date) is always read by Excel's type-sniffer as a date, so
    # force it to stay literal text with a leading quote-prefix -- same as
    # manually typing '2026-01-28 into a text-formatted log cell.
    $ws.Cells.Item($r, 1).Value = "'" + $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]

    if ($percentCol5) {
        $ws.Cells.Item($r, 5).Value = "'" + $entry[5]
    } else {
        $ws.Cells.Item($r, 5).Value = $entry[5]
    }

    $ws.Cells.Item($r, 6).Value = $entry[6]
}

$wsPir = $wb.Worksheets.Item("PIR")
foreach ($entry in $pirRows) {
    Write-LogRow $wsPir $entry $false
}

$wsHumidity = $wb.Worksheets.Item("Humidity")
foreach ($entry in $humidityRows) {
    Write-LogRow $wsHumidity $entry $true
}

$wsTemperature = $wb.Worksheets.Item("Temperature")
foreach ($entry in $temperatureRows) {
    Write-LogRow $wsTemperature $entry $false
}

Write-Host "Appended $($pirRows.Count) PIR rows, $($humidityRows.Count) Humidity rows, $($temperatureRows.Count) Temperature rows."
